$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")
$ws.Activate()

# Update "Responsible" (مسئول انجام) column for the three User-list rows
# from "سعید نوری" to the newly added "سعید قاسمی"
$ws.Range("E5").Value = "سعید قاسمی"
$ws.Range("E6").Value = "سعید قاسمی"
$ws.Range("E7").Value = "سعید قاسمی"

# Update the view: scroll back to top and move the selection to E8
$ws.Range("A1").Select() | Out-Null
$ws.Range("E8").Select() | Out-Null
